$d = $word.ActiveDocument

$replacements = @(
    @{old="741×3="; new="576×2="},
    @{old="299×7="; new="175×6="},
    @{old="556×4="; new="842×7="},
    @{old="758×4="; new="514×4="},
    @{old="519×8="; new="311×2="},
    @{old="633×4="; new="633×9="},
    @{old="941×6="; new="641×4="},
    @{old="541×3="; new="152×8="},
    @{old="716×2="; new="101×4="},
    @{old="834×5="; new="984×7="},
    @{old="992×5="; new="441×6="},
    @{old="250×9="; new="897×4="},
    @{old="539×2="; new="849×6="},
    @{old="513×8="; new="676×8="},
    @{old="568×2="; new="380×6="},
    @{old="575×8="; new="928×9="},
    @{old="622×9="; new="442×6="},
    @{old="500×4="; new="739×6="},
    @{old="594×6="; new="170×8="},
    @{old="949×8="; new="102×9="},
    @{old="350×5="; new="238×5="},
    @{old="937×7="; new="840×4="},
    @{old="966×7="; new="900×4="},
    @{old="425×5="; new="567×4="},
    @{old="857×2="; new="563×7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
